# Apply the "Aangekondigd" sheet commit:
#   1. On "Alle initiatieven", the running index in column A (rows 2-21)
#      is bumped by one (0..19 -> 1..20).
#   2. A new worksheet "Aangekondigd" is inserted right after it, containing
#      the very same header/data (now showing the updated 1..20 numbering),
#      with the data rows (B2:E21) given a new solid-fill highlight style.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- 1. bump column A on the existing sheet -------------------------------
for ($r = 2; $r -le 21; $r++) {
    $cell = $ws1.Cells.Item($r, 1)
    $cell.Value = [double]$cell.Value2 + 1
}

# --- 2. insert the new sheet right after "Alle initiatieven" -------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Aangekondigd"
$ws2.Outline.SummaryRow    = 1
$ws2.Outline.SummaryColumn = 1

# Clone the formatting (fonts/borders/alignment) of the header row and the
# numbering column from sheet1 onto the same cells of the new sheet.
$ws1.Range("A1:E21").Copy()
$ws2.Range("A1:E21").PasteSpecial(-4122)   # xlPasteFormats

# Copy over the (now renumbered) values/text cell by cell.
for ($r = 1; $r -le 21; $r++) {
    for ($c = 1; $c -le 5; $c++) {
        $src = $ws1.Cells.Item($r, $c)
        $dst = $ws2.Cells.Item($r, $c)
        $dst.Value = $src.Value2
    }
}

# Highlight the data cells (B2:E21) on the new sheet with a solid fill -
# this is the extra style/fill introduced by the commit.
$ws2.Range("B2:E21").Interior.Pattern = 1   # xlSolid

# Keep the original sheet as the active/selected one.
$ws1.Activate()
